$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 14723.533
$ws.Range("J19").Value = 20415.8
$ws.Range("L19").Value = 20415.8
$ws.Range("N19").Value = -20765.8
$ws.Range("H33").Value = 68399.17999999999
$ws.Range("I33").Value = 75225.2
$ws.Range("J33").Value = 139
$ws.Range("K33").Value = 75225.2
$ws.Range("L33").Value = 139
$ws.Range("M33").Value = -74996.2
$ws.Range("N33").Value = -597
$ws.Range("H39").Value = 147.15
$ws.Range("I39").Value = 140.93333
$ws.Range("K39").Value = 422.79999
$ws.Range("M39").Value = -126.79999
$ws.Range("H86").Value = 8343153.5
$ws.Range("I86").Value = 6449.3125
$ws.Range("K86").Value = 6449.3125
$ws.Range("M86").Value = -5326.3125
$ws.Range("H89").Value = 8343153.5
$ws.Range("I89").Value = 6449.3125
$ws.Range("K89").Value = 32246.5625
$ws.Range("M89").Value = -26630.5625
$ws.Range("H116").Value = 38968280
$ws.Range("I116").Value = 31378314
$ws.Range("K116").Value = 31378314
$ws.Range("M116").Value = -31374872
$ws.Range("H125").Value = 956.2857
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()
$ws.Range("H132").Value = 3357.0334
$ws.Range("I132").Value = 2769
$ws.Range("K132").Value = 8307
$ws.Range("M132").Value = -5777

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2243.4285
$ws.Range("I74").Value = 2107.647
$ws.Range("J74").Value = 2820.5
$ws.Range("K74").Value = 2107.647
$ws.Range("L74").Value = 2820.5
$ws.Range("M74").Value = -1233.647
$ws.Range("N74").Value = -4568.5
$ws.Range("H77").Value = 2243.4285
$ws.Range("I77").Value = 2107.647
$ws.Range("J77").Value = 2820.5
$ws.Range("K77").Value = 10538.235
$ws.Range("L77").Value = 14102.5
$ws.Range("M77").Value = -6170.235000000001
$ws.Range("N77").Value = -22838.5
$ws.Range("H88").Value = 12822871
$ws.Range("J88").Value = 2849.7778
$ws.Range("L88").Value = 2849.7778
$ws.Range("N88").Value = -3661.7778
$ws.Range("H91").Value = 12822871
$ws.Range("J91").Value = 2849.7778
$ws.Range("L91").Value = 2849.7778
$ws.Range("N91").Value = -5657.7778
$ws.Range("H102").Value = 904960
$ws.Range("I102").Value = 1251188.8
$ws.Range("K102").Value = 1251188.8
$ws.Range("M102").Value = -1249566.8
$ws.Range("H132").Value = 30305200
$ws.Range("I132").Value = 34484800
$ws.Range("J132").Value = 3098.5
$ws.Range("K132").Value = 103454400
$ws.Range("L132").Value = 9295.5
$ws.Range("M132").Value = -103451870
$ws.Range("N132").Value = -14355.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H102").Value = 0
$ws.Range("I102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("M102").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3149.6965
$ws.Range("I31").Value = 2292.9333
$ws.Range("J31").Value = 3463.1462
$ws.Range("K31").Value = 2292.9333
$ws.Range("L31").Value = 3463.1462
$ws.Range("M31").Value = -1997.9333
$ws.Range("N31").Value = -4053.1462
$ws.Range("H34").Value = 3149.6965
$ws.Range("I34").Value = 2292.9333
$ws.Range("J34").Value = 3463.1462
$ws.Range("K34").Value = 2292.9333
$ws.Range("L34").Value = 3463.1462
$ws.Range("M34").Value = -2090.9333
$ws.Range("N34").Value = -3867.1462
$ws.Range("H58").Value = 2036.3429
$ws.Range("I58").Value = 1061.2727
$ws.Range("K58").Value = 1061.2727
$ws.Range("M58").Value = -858.2727
$ws.Range("H94").Value = 660.94446
$ws.Range("I94").Value = 432.2
$ws.Range("J94").Value = 1804.6666
$ws.Range("K94").Value = 432.2
$ws.Range("L94").Value = 1804.6666
$ws.Range("M94").Value = 18.80000000000001
$ws.Range("N94").Value = -2706.6666
$ws.Range("H134").Value = 2226.28
$ws.Range("I134").Value = 1313.6111
$ws.Range("J134").Value = 4573.143
$ws.Range("K134").Value = 3940.8333
$ws.Range("L134").Value = 13719.429
$ws.Range("M134").Value = -1405.8333
$ws.Range("N134").Value = -18789.429
$ws.Range("H136").Value = 2036.3429
$ws.Range("I136").Value = 1061.2727
$ws.Range("K136").Value = 3183.8181
$ws.Range("M136").Value = -633.8181

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 856.5714
$ws.Range("I5").Value = 409.8
$ws.Range("J5").Value = 1104.7778
$ws.Range("K5").Value = 1229.4
$ws.Range("L5").Value = 3314.3334
$ws.Range("M5").Value = -1117.4
$ws.Range("N5").Value = -3538.3334
$ws.Range("H107").Value = 711.46155
$ws.Range("I107").Value = 891.125
$ws.Range("J107").Value = 424
$ws.Range("K107").Value = 2673.375
$ws.Range("L107").Value = 1272
$ws.Range("M107").Value = -753.375
$ws.Range("N107").Value = -5112
$ws.Range("H132").Value = 2039.3334
$ws.Range("J132").Value = 2344.5454
$ws.Range("L132").Value = 21100.9086
$ws.Range("N132").Value = -26160.9086
$ws.Range("H135").Value = 856.5714
$ws.Range("I135").Value = 409.8
$ws.Range("J135").Value = 1104.7778
$ws.Range("K135").Value = 3688.2
$ws.Range("L135").Value = 9943.0002
$ws.Range("M135").Value = -1153.2
$ws.Range("N135").Value = -15013.0002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2384.4583
$ws.Range("I80").Value = 2743.8
$ws.Range("J80").Value = 2127.7856
$ws.Range("K80").Value = 2743.8
$ws.Range("L80").Value = 2127.7856
$ws.Range("M80").Value = -1745.8
$ws.Range("N80").Value = -4123.7856
$ws.Range("H83").Value = 2384.4583
$ws.Range("I83").Value = 2743.8
$ws.Range("J83").Value = 2127.7856
$ws.Range("K83").Value = 13719
$ws.Range("L83").Value = 10638.928
$ws.Range("M83").Value = -8727
$ws.Range("N83").Value = -20622.928
$ws.Range("H126").Value = 10507.333
$ws.Range("J126").Value = 3366.5
$ws.Range("L126").Value = 10099.5
$ws.Range("N126").Value = -15039.5
$ws.Range("H132").Value = 4679.952
$ws.Range("I132").Value = 4406
$ws.Range("J132").Value = 5125.125
$ws.Range("K132").Value = 13218
$ws.Range("L132").Value = 15375.375
$ws.Range("M132").Value = -10688
$ws.Range("N132").Value = -20435.375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3999.7856
$ws.Range("I40").Value = 2999.625
$ws.Range("J40").Value = 5333.3335
$ws.Range("K40").Value = 2999.625
$ws.Range("L40").Value = 5333.3335
$ws.Range("M40").Value = -2863.625
$ws.Range("N40").Value = -5605.3335
$ws.Range("H68").Value = 29334
$ws.Range("I68").Value = 2999
$ws.Range("K68").Value = 2999
$ws.Range("M68").Value = -2250
$ws.Range("H71").Value = 29334
$ws.Range("I71").Value = 2999
$ws.Range("K71").Value = 14995
$ws.Range("M71").Value = -11251
$ws.Range("H133").Value = 107500
$ws.Range("J133").Value = 107500
$ws.Range("L133").Value = 107500
$ws.Range("N133").Value = -112560

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 101388.5
$ws.Range("J64").Value = 55000
$ws.Range("L64").Value = 55000
$ws.Range("N64").Value = -55496
$ws.Range("H67").Value = 101388.5
$ws.Range("J67").Value = 55000
$ws.Range("L67").Value = 55000
$ws.Range("N67").Value = -56716
$ws.Range("H122").Value = 2907.5386
$ws.Range("I122").Value = 2899.8
$ws.Range("J122").Value = 2933.3333
$ws.Range("K122").Value = 8699.400000000001
$ws.Range("L122").Value = 8799.999899999999
$ws.Range("M122").Value = -6249.400000000001
$ws.Range("N122").Value = -13699.9999
$ws.Range("H132").Value = 6872.125
$ws.Range("I132").Value = 8119.375
$ws.Range("J132").Value = 5624.875
$ws.Range("K132").Value = 24358.125
$ws.Range("L132").Value = 16874.625
$ws.Range("M132").Value = -21828.125
$ws.Range("N132").Value = -21934.625
$ws.Range("H136").Value = 3524.3157
$ws.Range("I136").Value = 1397.9333
$ws.Range("K136").Value = 4193.7999
$ws.Range("M136").Value = -1643.7999
